# Integrate observed overall FP votes into live simulations.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calc")

# --- ResolvePM 2PP conversion raw numbers (row 3) ---
$ws.Range("J3").Value = 30
$ws.Range("K3").Value = 37
$ws.Range("L3").Value = 13
$ws.Range("O3").Value = 6
$ws.Range("P3").Value = 6

# Update the 2PP conversion formula to normalise by the total first-preference vote
$ws.Range("Q3").Formula = "=(K3+L3*0.822+M3*0.348+N3*0.352+(O3+P3)*(0.507))*100/SUM(J3:P3)"

# --- Latest Newspoll (row 5) ---
$ws.Range("B5").Value = 51.9
$ws.Range("C5").Value = 52.38
$ws.Range("D5").Value = 53.39
$ws.Range("E5").Value = 46.85
$ws.Range("F5").Value = 51.64
$ws.Range("G5").Value = 53.69

# --- Latest Morgan (row 6) ---
$ws.Range("B6").Value = 53
$ws.Range("C6").Value = 52
$ws.Range("D6").Value = 57
$ws.Range("E6").Value = 47
$ws.Range("F6").Value = 54.5
$ws.Range("G6").Value = 49

# --- Second Morgan (row 7) ---
$ws.Range("C7").Value = 51.5
$ws.Range("D7").Value = 61
$ws.Range("E7").Value = 46.5
$ws.Range("F7").Value = 57.5
$ws.Range("G7").Value = 62.5

# --- Third Morgan (row 8) ---
$ws.Range("B8").Value = 55.5
$ws.Range("C8").Value = 56
$ws.Range("D8").Value = 63.5
$ws.Range("E8").Value = 43.5
$ws.Range("G8").Value = 62.5

# --- ResolvePM -> (row 9), new observed poll with live FP-derived 2PP ---
$ws.Range("B9").Value = 51.589898989898984
$ws.Range("C9").Value = 47.675510204081633
$ws.Range("D9").Value = 52.043999999999997
$ws.Range("E9").Value = 51.463000000000001
$ws.Range("H9").Value = 56.785858585858591

# --- Second ResolvePM -> (row 10), rolled down from the prior ResolvePM entry ---
$ws.Range("B10").Value = 53.885999999999996
$ws.Range("C10").Value = 57.556999999999995
$ws.Range("D10").Value = 53.905000000000001
$ws.Range("E10").Value = 47.676000000000002
$ws.Range("H10").Value = 53.067999999999998

# --- Third ResolvePM -> (row 11), rolled down from the prior Second ResolvePM entry ---
$ws.Range("B11").Value = 52.433000000000007
$ws.Range("C11").Value = 52.417999999999992
$ws.Range("D11").Value = 51.266999999999996
$ws.Range("E11").Value = 51.446999999999996
$ws.Range("H11").Value = 53.558999999999997

# Move the cell cursor, matching the author's final selection
$ws.Range("C6").Select()
